$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44319, 2, 34, 103.0834066034017),
    @(44320, 4, 30, 90.95594700300154),
    @(44321, 0, 29, 87.92408210290149)
)

$startRow = 245
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Range("A244").Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
